$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "changed black oxide parts to zinc allow[oy]"
# Solar panel mount screws: McMaster Carr part number (black-oxide) -> zinc-alloy equivalent
$ws.Range("D9").Value = "90128A106"

# Signpost attachment Bolts: McMaster Carr part number (black-oxide) -> zinc-alloy equivalent
$ws.Range("D10").Value = "90128A636"

# Reflect where the user left the cursor after making the edit
$ws.Range("D9").Select()
